# TC14 Canine Filter Breed-BullMastif: corrected the FilesTab Neo4j query
# (B4) to drop the redundant `File Type` and `Breed` return columns, and
# reselect that cell per the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newFilesQuery = "MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed  IN ['Bullmastiff']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS ``File Name``,
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``,
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newFilesQuery

# The cell shrank from 17 wrapped lines to 15 (two `coalesce` rows were
# removed), so the wrap-text row autofits from 246.5pt to 217.5pt
# (15 lines * 14.5pt default row height, matching the other rows on
# this sheet which are likewise plain multiples of the 14.5pt default).
$ws.Rows("4:4").RowHeight = 217.5

# Author ended up with B4 selected, scrolled so row 4 is in view.
$ws.Activate()
$ws.Range("B4").Select()
